$d = $word.ActiveDocument

# The bold run "DOCX, DOC, PDF, HTML, XPS, R" / hidden "_GoBack" bookmark /
# bold run "TF and TXT" are merged back into a single bold run whose text is
# "DOCX, DOC, PDF, HTML, XPS, RTF and TXT" (the bookmark, being inside the
# replaced range, is dropped as part of the replace).
$rng = $d.Content
$rng.Find.Execute("DOCX, DOC, PDF, HTML, XPS, RTF and TXT", $true, $false, $false, $false, $false, $true, 1, $false, `
                   "DOCX, DOC, PDF, HTML, XPS, RTF and TXT", 2) | Out-Null
